$d = $word.ActiveDocument

# The document has a single inline picture (the "Gross Plot Ratio" figure).
# Replace it with a hyperlink run whose visible text is the image's URL.
$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F01_Gross_Plot_Ratio.jpg?h=100%25&w=100%25"

$shape = $d.InlineShapes.Item(1)
$r = $shape.Range

# Remove the picture, leaving a collapsed range where it used to be.
$shape.Delete()

# Insert a hyperlink at that location, pointing at the image URL and
# displaying the URL itself as the link text.
$d.Hyperlinks.Add($r, $url, $null, $null, $url)
